$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("第二学年")

# Row 10: 支出 400 生活费 on 2018-09-07 (serial 43350)
$ws.Range("B10").Value = 7
$ws.Range("C10").Value = "支出"
$ws.Range("D10").Value = 400
$ws.Range("E10").Value = 43350
$ws.Range("F10").Value = "生活费"
$ws.Range("G10").Value = "生活费"

# Row 11: 支出 390 (=140+250) 其它 on 2018-09-14 (serial 43357)
$ws.Range("B11").Value = 8
$ws.Range("C11").Value = "支出"
$ws.Range("D11").Formula = "=140+250"
$ws.Range("E11").Value = 43357
$ws.Range("F11").Value = "其它"
$ws.Range("G11").Value = "补交住宿费200+车费50+水费140"

$ws.Range("C10:C11").HorizontalAlignment = -4108

$ws.Range("J4").Calculate()
$ws.Range("K4").Calculate()
$ws.Range("J10").Calculate()
$ws.Range("J11").Calculate()

$ws.Range("I15").Select()
